$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.0253360835587989
$ws.Range("C2").Value = 1.135698984415373
$ws.Range("D2").Value = 7.612000487694052
$ws.Range("E2").Value = 2.758985409112207
$ws.Range("F2").Value = 2.78578550780152
$ws.Range("G2").Value = 52
$ws.Range("B3").Value = 0.08600705238178873
$ws.Range("C3").Value = 1.123115432295725
$ws.Range("D3").Value = 5.609313565901017
$ws.Range("E3").Value = 2.368398945680608
$ws.Range("F3").Value = 2.390387976856151
$ws.Range("G3").Value = 51
$ws.Range("B4").Value = 0.0363203761572407
$ws.Range("C4").Value = 0.9693363125750792
$ws.Range("D4").Value = 4.390696811956634
$ws.Range("E4").Value = 2.095398962478658
$ws.Range("F4").Value = 2.116354596420031
$ws.Range("G4").Value = 50
$ws.Range("B5").Value = 0.1260652424010658
$ws.Range("C5").Value = 1.03221972966461
$ws.Range("D5").Value = 4.978658809889423
$ws.Range("E5").Value = 2.231290839377382
$ws.Range("F5").Value = 2.25081259337064
$ws.Range("G5").Value = 49
$ws.Range("B6").Value = 0.08428106382731648
$ws.Range("C6").Value = 1.007929901751503
$ws.Range("D6").Value = 5.095203452855081
$ws.Range("E6").Value = 2.257255734925726
$ws.Range("F6").Value = 2.27955206241149
$ws.Range("G6").Value = 48
$ws.Range("B7").Value = 0.09619512034022865
$ws.Range("C7").Value = 1.131065357075751
$ws.Range("D7").Value = 5.749733647272278
$ws.Range("E7").Value = 2.397860222630226
$ws.Range("F7").Value = 2.42991643159901
$ws.Range("G7").Value = 36
$ws.Range("B8").Value = 0.1381149580396931
$ws.Range("C8").Value = 1.187418322474588
$ws.Range("D8").Value = 5.903068426619335
$ws.Range("E8").Value = 2.4296231038207
$ws.Range("F8").Value = 2.461107737049912
$ws.Range("G8").Value = 35
$ws.Range("B9").Value = 0.1382699113425551
$ws.Range("C9").Value = 1.562389918879535
$ws.Range("D9").Value = 9.831370974918348
$ws.Range("E9").Value = 3.135501710240061
$ws.Range("F9").Value = 3.223266002762909
$ws.Range("G9").Value = 18
$ws.Range("B10").Value = -0.7150897714758779
$ws.Range("C10").Value = 1.341686922679743
$ws.Range("D10").Value = 6.639370992174167
$ws.Range("E10").Value = 2.576697691265735
$ws.Range("F10").Value = 2.596308797503721
$ws.Range("G10").Value = 11
$ws.Range("B11").Value = 0.2236099196487949
$ws.Range("C11").Value = 0.5232375926460167
$ws.Range("D11").Value = 0.3429526541974182
$ws.Range("E11").Value = 0.5856215964233373
$ws.Range("F11").Value = 0.6051355819484564
$ws.Range("G11").Value = 5
